# ICT.xlsx - "update in attendance feature"
# Re-enter the attendance table (rows 2-7), fix a couple of per-cell style
# glitches caused by the data reshuffle, resize header rows, and move the
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Style fixes: a few cells need to switch to a style that already
#    exists elsewhere on the sheet. Copy number-format/font/alignment
#    from a stable donor cell so we reuse the existing style slot
#    instead of minting a new one. (Done one cell at a time -
#    PasteSpecial onto a multi-area "A1,B2" range only hits the first
#    area in this host.)
# ---------------------------------------------------------------------
foreach ($addr in @("F4", "F5", "F6", "F7", "G7", "H7", "H6")) {
    $ws.Range("F2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

foreach ($addr in @("A3", "I5", "J5", "I6", "J6", "I7", "J7")) {
    $ws.Range("B2").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2) Row heights / header row resize
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 11.25
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 17.25
$ws.Rows.Item(7).RowHeight = 34.5

# ---------------------------------------------------------------------
# 3) Re-enter the table contents
# ---------------------------------------------------------------------

# Header row
$ws.Range("A1").Value = "Course"
$ws.Range("B1").Value = "Department"
$ws.Range("C1").Value = "To Date"
$ws.Range("D1").Value = "From Date"
$ws.Range("E1").Value = "Class"
$ws.Range("F1").Value = "Status"
$ws.Range("G1").Value = "First Name"
$ws.Range("H1").Value = "Last Name"
$ws.Range("I1").Value = "Teacher First Name"
$ws.Range("J1").Value = "Teacher Last Name"

# Row 2
$ws.Range("A2").Value = "Computer Network"
$ws.Range("B2").Value = "Cs"
$ws.Range("C2").Value = 45541
$ws.Range("D2").Value = 45505
$ws.Range("E2").Value = "F-19"
$ws.Range("F2").Value = "PRESENT"
$ws.Range("G2").Value = "shaista"
$ws.Range("H2").Value = "tabbasum"
$ws.Range("I2").Value = "atiqa"
$ws.Range("J2").Value = "tabbasum"

# Row 3
$ws.Range("A3").Value = "arts&crafts"
$ws.Range("B3").Value = "Finanace"
$ws.Range("C3").Value = 45528
$ws.Range("D3").Value = 45527
$ws.Range("E3").Value = "G-11"
$ws.Range("F3").Value = "ABSENT"
$ws.Range("G3").Value = "marwa"
$ws.Range("H3").Value = "jabeen"
$ws.Range("I3").Value = "anfal"
$ws.Range("J3").Value = "Tabbasum"

# Row 4
$ws.Range("A4").Value = "urdu"
$ws.Range("B4").Value = "URDU"
$ws.Range("C4").Value = 45530
$ws.Range("D4").Value = 45505
$ws.Range("E4").Value = "G-11"
$ws.Range("F4").Value = "ABSENT"
$ws.Range("G4").Value = "alisha"
$ws.Range("H4").Value = "Fatima"
$ws.Range("I4").Value = "maham "
$ws.Range("J4").Value = "masoom"

# Row 5 (was blank)
$ws.Range("A5").Value = "ICT"
$ws.Range("B5").Value = "CS"
$ws.Range("C5").Value = 45528
$ws.Range("D5").Value = 45505
$ws.Range("E5").Value = "F-19"
$ws.Range("F5").Value = "PRESENT"
$ws.Range("G5").Value = "husnain"
$ws.Range("H5").Value = "ahmad"
$ws.Range("I5").Value = "ibtasam"
$ws.Range("J5").Value = "ilahi"

# Row 6 (was blank)
$ws.Range("A6").Value = "ICT"
$ws.Range("B6").Value = "math"
$ws.Range("C6").Value = 45541
$ws.Range("D6").Value = 45505
$ws.Range("E6").Value = "F-19"
$ws.Range("F6").Value = "PRESENT"
$ws.Range("G6").Value = "shaista"
$ws.Range("H6").Value = "Tabbasum"
$ws.Range("I6").Value = "atiqa"
$ws.Range("J6").Value = "tabbasum"

# Row 7 (was blank)
$ws.Range("A7").Value = "Computer Network"
$ws.Range("B7").Value = "CS"
$ws.Range("C7").Value = 45528
$ws.Range("D7").Value = 45505
$ws.Range("E7").Value = "F-19"
$ws.Range("F7").Value = "PRESENT"
$ws.Range("G7").Value = "Atiqa"
$ws.Range("H7").Value = "Tabbasum"
$ws.Range("I7").Value = "Lalain"
$ws.Range("J7").Value = "Fatima"

# ---------------------------------------------------------------------
# 4) Move the active selection
# ---------------------------------------------------------------------
$ws.Range("K10").Select() | Out-Null
